# "update po dan pr" - rename the J1/K1 headers on the PO import format sheet
# from "Name" / "Product name" to "Part name" / "Part number".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J1").Value = "Part name"
$ws.Range("K1").Value = "Part number"

# Columns J and K are best-fit to their header text; widen them to fit the
# new, longer headers (values chosen so the resulting stored column width
# matches the new best-fit widths as closely as possible).
$ws.Columns.Item(10).ColumnWidth = 12
$ws.Columns.Item(11).ColumnWidth = 14.166666666666666

# Leave the cursor where the author left it after editing.
$ws.Range("B7").Select()
